# "diem danh ngay 24" - add attendance column for 24-Mar-2012 (serial 40992)
# by mirroring the existing 23-Mar-2012 column (G) into a new column (H).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column G (dates/marks for 23-Mar) holds the layout/styling we want to
# replicate for the new day column H: copy values + formatting in one shot.
$ws.Range("G1:G7").Copy($ws.Range("H1:H7")) | Out-Null

# The header date must advance by one day (24-Mar-2012), not repeat G1's date.
$ws.Range("H1").Value2 = 40992

# Give column H the same display width as the other attendance columns.
$ws.Columns.Item(8).ColumnWidth = $ws.Columns.Item(7).ColumnWidth

# Move the active selection to the newly-filled-in cell, matching the
# author's cursor position after making the edit.
$ws.Range("H7").Select() | Out-Null
